$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L2").Value = 1074.68
$wsGrupo.Range("M2").Value = 4971.2
$wsGrupo.Range("L37").Value = "4 de 35"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F2").Value = 6045.88
$wsMensual.Range("F37").Value = 37836.15

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 11 - PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 2255.37
$wsCumpl.Range("E11").Value = 666.85458185274
$wsCumpl.Range("F11").Value = 0.7717989965610573

# Row 12 - PORCELANATO
$wsCumpl.Range("D12").Value = 33617.57
$wsCumpl.Range("E12").Value = -11916.3
$wsCumpl.Range("F12").Value = 1.549106112222925

# Row 14 - TOTAL
$wsCumpl.Range("D14").Value = 39238.05
$wsCumpl.Range("E14").Value = -2652.482762818175
$wsCumpl.Range("F14").Value = 1.072500796437631
